$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the .java file names to reflect the new package structure
$ws.Range("B2").Value = ".\compositor\Compositor.java"
$ws.Range("B3").Value = ".\compositor\TeXCompositor.java"
$ws.Range("B4").Value = ".\compositor\SimpleCompositor.java"
$ws.Range("B6").Value = ".\compositor\ArrayCompositor.java"
$ws.Range("B5").Value = ".\composition\Composition.java"

# Update the relationship cells to reflect the new "Import" relationship
$ws.Range("C3").Value = "Import,Implement"
$ws.Range("C4").Value = "Import,Implement"
$ws.Range("C5").Value = "Implement"
$ws.Range("C6").Value = "Contain,Import,Call,Use:2,Parameter"
$ws.Range("D7").Value = "Create,Import,Call"
$ws.Range("E7").Value = "Create,Import,Call"
$ws.Range("F7").Value = "Create,Import,Call"
$ws.Range("G7").Value = "Create:2,Import,Call:8,Contain:2"
